$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the three header cells (A1/B1/C1) to drop the extra ".1" in the
# indicator code (4.c.1.1 -> 4.c.1), matching the new shared-string text.
$ws.Range("A1").Value = "4.с.1 Билим берүү мекемелерде диплом берилгем мугалимдердин үлүшү"
$ws.Range("B1").Value = "4.c.1 Доля дипломированных учителей в образовательных учереждениях"
$ws.Range("C1").Value = "4.c.1 Proportion of certified teachers in educational institutions"

# Update the selected cell/active cell shown in the sheet view from N11 to C11.
$ws.Range("C11").Select()
